$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1579.1765
$ws.Range("I112").Value = 560
$ws.Range("J112").Value = 1677.8064
$ws.Range("K112").Value = 1680
$ws.Range("L112").Value = 5033.4192
$ws.Range("M112").Value = -572
$ws.Range("N112").Value = -7249.4192
$ws.Range("H132").Value = 28822.117
$ws.Range("I132").Value = 3080.6086
$ws.Range("K132").Value = 9241.825800000001
$ws.Range("M132").Value = -6711.825800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10854.884
$ws.Range("I32").Value = 9909.017
$ws.Range("K32").Value = 9909.017
$ws.Range("M32").Value = -9622.017
$ws.Range("H122").Value = 2061.8667
$ws.Range("I122").Value = 1757.1111
$ws.Range("J122").Value = 2519
$ws.Range("K122").Value = 5271.3333
$ws.Range("L122").Value = 7557
$ws.Range("M122").Value = -2821.3333
$ws.Range("N122").Value = -12457
$ws.Range("H123").Value = 38000
$ws.Range("J123").Value = 38000
$ws.Range("L123").Value = 38000
$ws.Range("N123").Value = -47800
$ws.Range("H132").Value = 23812362
$ws.Range("I132").Value = 45456920
$ws.Range("J132").Value = 3349.4
$ws.Range("K132").Value = 136370760
$ws.Range("L132").Value = 10048.2
$ws.Range("M132").Value = -136368230
$ws.Range("N132").Value = -15108.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2077.476
$ws.Range("I3").Value = 2360.0667
$ws.Range("J3").Value = 1371
$ws.Range("K3").Value = 2360.0667
$ws.Range("L3").Value = 1371
$ws.Range("M3").Value = -2246.0667
$ws.Range("N3").Value = -1599
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H134").Value = 3290.6897
$ws.Range("I134").Value = 2825.65
$ws.Range("J134").Value = 4324.1113
$ws.Range("K134").Value = 8476.950000000001
$ws.Range("L134").Value = 12972.3339
$ws.Range("M134").Value = -5941.950000000001
$ws.Range("N134").Value = -18042.3339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 52999
$ws.Range("J116").Value = 52999
$ws.Range("L116").Value = 52999
$ws.Range("N116").Value = -62177
$ws.Range("H132").Value = 402572.06
$ws.Range("I132").Value = 1500.5
$ws.Range("J132").Value = 1561223.2
$ws.Range("K132").Value = 4501.5
$ws.Range("L132").Value = 4683669.6
$ws.Range("M132").Value = -1971.5
$ws.Range("N132").Value = -4688729.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1658.5
$ws.Range("I34").Value = 450
$ws.Range("J34").Value = 2061.3333
$ws.Range("K34").Value = 1350
$ws.Range("L34").Value = 6183.999899999999
$ws.Range("M34").Value = -1266
$ws.Range("N34").Value = -6351.999899999999
$ws.Range("H39").Value = 600
$ws.Range("J39").Value = 600
$ws.Range("L39").Value = 1800
$ws.Range("N39").Value = -2388
$ws.Range("H55").Value = 2000
$ws.Range("J55").Value = 2000
$ws.Range("L55").Value = 6000
$ws.Range("N55").Value = -6354
$ws.Range("H68").Value = 1291.8701
$ws.Range("I68").Value = 1021.9
$ws.Range("J68").Value = 1386.5964
$ws.Range("K68").Value = 3065.7
$ws.Range("L68").Value = 4159.789199999999
$ws.Range("M68").Value = -2254.7
$ws.Range("N68").Value = -5781.789199999999
$ws.Range("H71").Value = 1291.8701
$ws.Range("I71").Value = 1021.9
$ws.Range("J71").Value = 1386.5964
$ws.Range("K71").Value = 9197.1
$ws.Range("L71").Value = 12479.3676
$ws.Range("M71").Value = -5141.1
$ws.Range("N71").Value = -20591.3676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1445.1818
$ws.Range("I122").Value = 1316.1666
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 3948.4998
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -1498.4998
$ws.Range("N122").Value = -9700
$ws.Range("H126").Value = 7613.1
$ws.Range("J126").Value = 2347.2222
$ws.Range("L126").Value = 7041.6666
$ws.Range("N126").Value = -11981.6666
$ws.Range("H132").Value = 50005692
$ws.Range("J132").Value = 4662.7
$ws.Range("L132").Value = 13988.1
$ws.Range("N132").Value = -19048.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1965.7059
$ws.Range("I61").Value = 1779.1111
$ws.Range("J61").Value = 2175.625
$ws.Range("K61").Value = 1779.1111
$ws.Range("L61").Value = 2175.625
$ws.Range("M61").Value = -1577.1111
$ws.Range("N61").Value = -2579.625
$ws.Range("H106").Value = 31500
$ws.Range("J106").Value = 31500
$ws.Range("L106").Value = 31500
$ws.Range("N106").Value = -34024
$ws.Range("H113").Value = 1965.7059
$ws.Range("I113").Value = 1779.1111
$ws.Range("J113").Value = 2175.625
$ws.Range("K113").Value = 1779.1111
$ws.Range("L113").Value = 2175.625
$ws.Range("M113").Value = 390.8888999999999
$ws.Range("N113").Value = -6515.625
$ws.Range("H122").Value = 65013
$ws.Range("I122").Value = 85700.664
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 257101.992
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -254651.992
$ws.Range("N122").Value = -13750

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1526.8667
$ws.Range("I96").Value = 1700.2727
$ws.Range("J96").Value = 1050
$ws.Range("K96").Value = 1700.2727
$ws.Range("L96").Value = 1050
$ws.Range("M96").Value = -327.2727
$ws.Range("N96").Value = -3796
$ws.Range("H122").Value = 2198863.2
$ws.Range("I122").Value = 7144058
$ws.Range("J122").Value = 998.8889
$ws.Range("K122").Value = 21432174
$ws.Range("L122").Value = 2996.6667
$ws.Range("M122").Value = -21429724
$ws.Range("N122").Value = -7896.6667
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 1840050.1
$ws.Range("I126").Value = 2102558.8
$ws.Range("J126").Value = 2490
$ws.Range("K126").Value = 6307676.399999999
$ws.Range("L126").Value = 7470
$ws.Range("M126").Value = -6305206.399999999
$ws.Range("N126").Value = -12410
$ws.Range("H132").Value = 1280636
$ws.Range("I132").Value = 2175282.8
$ws.Range("J132").Value = 2568.9285
$ws.Range("K132").Value = 6525848.399999999
$ws.Range("L132").Value = 7706.7855
$ws.Range("M132").Value = -6523318.399999999
$ws.Range("N132").Value = -12766.7855
